$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. First paragraph: append a parenthetical note in red, split across
#    three runs, after padding the original sentence with two spaces.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Range = $d.Range($p1.Range.Start, $p1.Range.End)

$newPara1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t xml:space="preserve">This is a Microsoft word document.  </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>(This is a change – Ve</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>rsion for branch alternate</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="C00000"/></w:rPr><w:t>)</w:t></w:r>' +
  '</w:p>'

$p1Range.InsertXML($newPara1Xml)

# ---------------------------------------------------------------------
# 2. Append a new, empty, shaded paragraph after the last paragraph in
#    the document body.
# ---------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$newLastParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr>' +
  '</w:p>'

$endRange.InsertXML($newLastParaXml)

Write-Host "Edit applied successfully"
